$d = $word.ActiveDocument

$d.Content.Find.Execute("Ben Barrrrr", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ben Bar", 2)
